# Updated cryptos list on Sun Sep 22 04:34:15 UTC 2024 with GitHub Actions
# Applies the per-cell price / volume(1h) refresh, plus the PEPE <-> Bittensor
# rank swap (rows 33/34), exactly as captured by the source diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-coerced into a numeric value
# by Excel (losing a trailing zero / exact decimal text, or changing cell type
# from text to number) are first forced to Text format so the literal string
# is preserved verbatim, matching the source inline-string cells.

$ws.Range('D2').Value = '63.127.18'
$ws.Range('E2').Value = '  +0.37%  '
$ws.Range('D3').Value = '2.597.10'
$ws.Range('E3').Value = '  +2.41%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '583.29'
$ws.Range('E5').Value = '  +2.64%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '148.67'
$ws.Range('E6').Value = '  +2.32%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.598'
$ws.Range('E8').Value = '  +2.91%  '
$ws.Range('E9').Value = '  +3.47%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '5.67'
$ws.Range('E10').Value = '  +3.47%  '
$ws.Range('E11').Value = '  +0.16%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.354'
$ws.Range('E12').Value = '  +0.48%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '27.34'
$ws.Range('E13').Value = '  +0.83%  '
$ws.Range('D14').Value = '3.062.96'
$ws.Range('E14').Value = '  +2.49%  '
$ws.Range('D15').Value = '62.991.08'
$ws.Range('E15').Value = '  +0.29%  '
$ws.Range('E16').Value = '  +4.17%  '
$ws.Range('D17').Value = '2.593.48'
$ws.Range('E17').Value = '  +2.02%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '11.39'
$ws.Range('E18').Value = '  +1.29%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '344.48'
$ws.Range('E19').Value = '  +3.38%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '4.42'
$ws.Range('E20').Value = '  +2.66%  '
$ws.Range('E21').Value = '  +1.18%  '
$ws.Range('E22').Value = '  -0.10%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.70'
$ws.Range('E23').Value = '  -0.96%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '67.25'
$ws.Range('E24').Value = '  +3.08%  '
$ws.Range('D25').Value = '2.723.24'
$ws.Range('E25').Value = '  +2.50%  '
$ws.Range('E26').Value = '  +0.06%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '1.61'
$ws.Range('E27').Value = '  +1.33%  '
$ws.Range('E28').Value = '  -0.21%  '
$ws.Range('E29').Value = '  +1.41%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '7.86'
$ws.Range('E30').Value = '  +8.52%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.45'
$ws.Range('E31').Value = '  +0.54%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.94'
$ws.Range('E32').Value = '  +5.65%  '
$ws.Range('B33').Value = 'Bittensor'
$ws.Range('C33').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '471.43'
$ws.Range('E33').Value = '  +18.20%  '
$ws.Range('B34').Value = 'PEPE'
$ws.Range('C34').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D34').Value = '0.0₃0825'
$ws.Range('E34').Value = '  +2.51%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '176.66'
$ws.Range('E35').Value = '  +0.03%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.62'
$ws.Range('E36').Value = '  +5.70%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.406'
$ws.Range('E37').Value = '  +2.01%  '
$ws.Range('E38').Value = '  -0.01%  '
$ws.Range('E39').Value = '  +0.96%  '
$ws.Range('E40').Value = '  +6.87%  '
$ws.Range('E41').Value = '  +0.02%  '
$ws.Range('E42').Value = '  -1.75%  '
$ws.Range('E43').Value = '  +5.81%  '
$ws.Range('E44').Value = '  +2.53%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.636'
$ws.Range('E45').Value = '  +6.64%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '21.23'
$ws.Range('E46').Value = '  +2.98%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0550'
$ws.Range('E47').Value = '  +4.05%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.0974'
$ws.Range('E48').Value = '  +1.32%  '
$ws.Range('E49').Value = '  +1.13%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '18.61'
$ws.Range('E50').Value = '  +2.81%  '
$ws.Range('E51').Value = '  +3.65%  '
